$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(394, 6).Value = 110.53
$ws.Cells.Item(395, 6).Value = 721.79
$ws.Cells.Item(396, 6).Value = 716.18
$ws.Cells.Item(397, 6).Value = 667.05
$ws.Cells.Item(398, 6).Value = 634.12
$ws.Cells.Item(399, 6).Value = 126.58
$ws.Cells.Item(400, 6).Value = 40.99
$ws.Cells.Item(401, 6).Value = 795.05
$ws.Cells.Item(402, 6).Value = 715.63
$ws.Cells.Item(403, 6).Value = 711.3
$ws.Cells.Item(404, 6).Value = 672.28
$ws.Cells.Item(405, 6).Value = 625.4
$ws.Cells.Item(406, 6).Value = 126.91
$ws.Cells.Item(407, 6).Value = 41.07
$ws.Cells.Item(408, 6).Value = 114.56
$ws.Cells.Item(409, 6).Value = 705.21
$ws.Cells.Item(410, 6).Value = 696.63
$ws.Cells.Item(411, 6).Value = 663.53
$ws.Cells.Item(412, 6).Value = 639.64
$ws.Cells.Item(413, 6).Value = 126.89
$ws.Cells.Item(414, 6).Value = 72.86
$ws.Cells.Item(415, 6).Value = 817.05
$ws.Cells.Item(416, 6).Value = 705.08
$ws.Cells.Item(417, 6).Value = 744.02
$ws.Cells.Item(418, 6).Value = 651.67999999999995
$ws.Cells.Item(419, 6).Value = 604.9
$ws.Cells.Item(420, 6).Value = 133.04
$ws.Cells.Item(421, 6).Value = 38.25
$ws.Cells.Item(422, 6).Value = 832.59
$ws.Cells.Item(423, 6).Value = 727.77
$ws.Cells.Item(424, 6).Value = 736.04
$ws.Cells.Item(425, 6).Value = 678.65
$ws.Cells.Item(426, 6).Value = 648.67999999999995
$ws.Cells.Item(427, 6).Value = 125.79
$ws.Cells.Item(428, 6).Value = 42.5
$ws.Cells.Item(429, 6).Value = 801.36
$ws.Cells.Item(430, 6).Value = 738.12
$ws.Cells.Item(431, 6).Value = 733.96
$ws.Cells.Item(432, 6).Value = 676.95
$ws.Cells.Item(433, 6).Value = 734.35
$ws.Cells.Item(434, 6).Value = 125.84
$ws.Cells.Item(435, 6).Value = 41.45
$ws.Cells.Item(436, 6).Value = 1074.45
$ws.Cells.Item(437, 6).Value = 784.76
$ws.Cells.Item(438, 6).Value = 729.95
$ws.Cells.Item(439, 6).Value = 691.33
$ws.Cells.Item(440, 6).Value = 637.4
$ws.Cells.Item(441, 6).Value = 119.65
$ws.Cells.Item(442, 6).Value = 43.44
$ws.Cells.Item(443, 6).Value = 837.28
$ws.Cells.Item(444, 6).Value = 740.52
$ws.Cells.Item(445, 6).Value = 730.82
$ws.Cells.Item(446, 6).Value = 674.73
$ws.Cells.Item(447, 6).Value = 607.74
$ws.Cells.Item(448, 6).Value = 122.41
$ws.Cells.Item(449, 6).Value = 43.48
$ws.Cells.Item(450, 6).Value = 799.27
$ws.Cells.Item(451, 6).Value = 741.39
$ws.Cells.Item(452, 6).Value = 732.08
$ws.Cells.Item(453, 6).Value = 641.66
$ws.Cells.Item(454, 6).Value = 669.06
$ws.Cells.Item(455, 6).Value = 123.2
$ws.Cells.Item(456, 6).Value = 40.82
$ws.Cells.Item(457, 6).Value = 843.58
$ws.Cells.Item(458, 6).Value = 746.85
$ws.Cells.Item(459, 6).Value = 699.4
$ws.Cells.Item(460, 6).Value = 726.7
$ws.Cells.Item(461, 6).Value = 623.28
$ws.Cells.Item(462, 6).Value = 124.85
$ws.Cells.Item(463, 6).Value = 42.22
$ws.Cells.Item(464, 6).Value = 1004.77
$ws.Cells.Item(465, 6).Value = 739.78
$ws.Cells.Item(466, 6).Value = 727.51
$ws.Cells.Item(467, 6).Value = 726.1
$ws.Cells.Item(468, 6).Value = 623.82000000000005
$ws.Cells.Item(469, 6).Value = 121.45
$ws.Cells.Item(470, 6).Value = 41.89
$ws.Cells.Item(471, 6).Value = 1004.06
$ws.Cells.Item(472, 6).Value = 741.36
$ws.Cells.Item(473, 6).Value = 750.2
$ws.Cells.Item(474, 6).Value = 719.51
$ws.Cells.Item(475, 6).Value = 639.16
$ws.Cells.Item(476, 6).Value = 121.41
$ws.Cells.Item(477, 6).Value = 43.55
$ws.Cells.Item(478, 6).Value = 1012.34
$ws.Cells.Item(479, 6).Value = 747.12
$ws.Cells.Item(480, 6).Value = 752.06
$ws.Cells.Item(481, 6).Value = 728.58
$ws.Cells.Item(482, 6).Value = 726.81
$ws.Cells.Item(483, 6).Value = 123.92
$ws.Cells.Item(484, 6).Value = 35.72
$ws.Cells.Item(485, 6).Value = 835.69
$ws.Cells.Item(486, 6).Value = 785.42
$ws.Cells.Item(487, 6).Value = 752.35
$ws.Cells.Item(488, 6).Value = 753.58
$ws.Cells.Item(489, 6).Value = 648.86

$ws.Range("F394:F489").Select()
